$wb = $excel.ActiveWorkbook

# Activate the "Inputs and Outputs" sheet and update its title cell.
$ws = $wb.Worksheets.Item("Inputs and Outputs")
$ws.Activate()
$ws.Range("A1").Value = "Results Summary and Inputs"
$ws.Range("A2").Select()
